$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.292562961578369
$ws.Range("B1").Value = 2.831019401550293
$ws.Range("C1").Value = 2.27980899810791
$ws.Range("D1").Value = 2.160499095916748
$ws.Range("E1").Value = 1.876410126686096
